$wb = $excel.ActiveWorkbook

# --- Sheet "Components & input parameter": rename header and convert in/out text to boolean 0/1 ---
$ws2 = $wb.Worksheets.Item("Components & input parameter")

$ws2.Range("B1").Value = "Inside"

$ws2.Range("B2").Value = 0
$ws2.Range("B3").Value = 1
$ws2.Range("B4").Value = 1
$ws2.Range("B5").Value = 0

# Make this sheet the active / selected sheet, matching the new selection.
$ws2.Activate()
$ws2.Range("B6").Select()

$wb.Save()
